$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: login/config row updated after integrating DB + new automation script
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "http://testpanvelmc.ptaxcollection.com:8080/Pages/Login.aspx"
$ws.Range("D2").Value = "KM"
$ws.Range("E2").Value = 6

# Rows 10-13: node/property values updated from BMC/old codes to KH/new codes
$ws.Range("A10").Value = "KH"
$ws.Range("C10").Value = "1-101"

$ws.Range("A11").Value = "KH"
$ws.Range("C11").Value = "1-14"

$ws.Range("A12").Value = "KH"
$ws.Range("C12").Value = "1-26"

$ws.Range("A13").Value = "KH"
$ws.Range("C13").Value = "1-28"

# Update active selection to match the authored workbook state
$ws.Range("C14").Select()
